# Automated monthly rollover update.
# Sheet "VENTAS POR GRUPO" (current-month per-product-group sales): zero out
# the cells that represented the just-closed month's sales, and refresh the
# "X de 31" completion counters for the affected columns.
# Sheet "VENTA MENSUAL" (rolling 4-month view): shift abril/mayo/junio/julio
# -> mayo/junio/julio/agosto, i.e. each row's C:F values slide one column to
# the left and the new rightmost month (agosto) starts at 0.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Narrow column J (LED) from width 11 to width 9 (COM ColumnWidth has a
# constant +0.8333333333333334 offset vs. the raw OOXML column width).
$ws1.Columns.Item(10).ColumnWidth = 8.166666666666666

# Cells whose value goes to 0.
$ws1ZeroCells = @(
    "H2", "M2",
    "D4", "L4", "M4",
    "M11", "O11", "P11",
    "M12",
    "M15",
    "D19", "E19", "J19", "M19", "O19", "P19",
    "D22", "M22",
    "M23",
    "M29"
)
foreach ($ref in $ws1ZeroCells) {
    $ws1.Range($ref).Value = 0
}

# Row 33 "X de 31" completion counters for the columns touched above.
$ws1Row33Cols = @("D", "E", "H", "J", "L", "M", "O", "P")
foreach ($col in $ws1Row33Cols) {
    $ws1.Range("$col" + "33").Value = "0 de 31"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths for C, E, F change (D and G stay the same).
$ws2.Columns.Item(3).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(5).ColumnWidth = 12.166666666666666
$ws2.Columns.Item(6).ColumnWidth = 11.166666666666666

# Month headers shift by one: abril/mayo/junio/julio -> mayo/junio/julio/agosto
$ws2.Range("C1").Value = "mayo"
$ws2.Range("D1").Value = "junio"
$ws2.Range("E1").Value = "julio"
$ws2.Range("F1").Value = "agosto"

# Each data row (and the totals row 33) shifts C:F left by one column; the
# new column F (agosto) starts at 0. Values below are the *old* D, E, F for
# that row, which become the *new* C, D, E.
$ws2Shift = @(
    @(2, 3163.97, 2862.27, 4360.63),
    @(3, 56.32, 0, 0),
    @(4, 285.12, 285.12, 1190.78),
    @(5, 0, 0, 0),
    @(6, 0, 10.28, 10.28),
    @(7, 0, 0, 0),
    @(8, 1565.15, 3068.54, 0),
    @(9, 0, 0, 0),
    @(10, 0, 0, 0),
    @(11, 0, 0, -86.23),
    @(12, 32.36, 3997, 112.01),
    @(13, 0, 0, 0),
    @(14, 0, 0, 0),
    @(15, 2511.66, 6533.02, 374.03),
    @(16, 0, 0, 0),
    @(17, 0, 0, 0),
    @(18, 0, 0, 0),
    @(19, 1509.73, 1974.55, 5850.44),
    @(20, 0, 0, 0),
    @(21, 0, 0, 0),
    @(22, 0, 1153.42, 612.28),
    @(23, 62.44, 0, 128.3),
    @(24, 0, 0, 0),
    @(25, 0, 0, 0),
    @(26, 222.19, 0, 0),
    @(27, 12043.47, 11532.42, 0),
    @(28, 0, 0, 0),
    @(29, 0, 0, 10.76),
    @(30, 0, 0, 0),
    @(31, 0, 0, 0),
    @(32, 0, 41.6, 0),
    @(33, 21452.41, 31458.22, 12563.28)
)

foreach ($row in $ws2Shift) {
    $r = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $ws2.Cells.Item($r, 6).Value = 0
}
